# Fix: Fixed a bug in the import process.
#
# The "properties_list" translation table on Sheet1 was missing the
# "wizard.new.import.database" entry. Insert it (Key / English /
# Japanese = "wizard.new.import.database" / "Database" / "データベース")
# in its correct alphabetically-sorted position - immediately above the
# existing "wizard.new.import.driver" row, which is row 263. All rows
# from 263 downward shift down by one (263 -> 264, ..., 274 -> 275) and
# the sheet dimension / table / autofilter ranges grow from C274 to
# C275 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$insertRow = 263

# Copy the row that is about to be pushed down (current row 263) and
# "insert copied cells" in its own place: this shifts it (and every
# row below it) down by one AND clones its formatting/content into the
# freshly-opened row 263, so the new row starts out styled exactly
# like its neighbours.
$ws.Rows.Item($insertRow).Copy()
$ws.Rows.Item($insertRow).Insert()
$excel.CutCopyMode = 0

# Overwrite the newly inserted row 263 with the new translation entry.
$ws.Cells.Item($insertRow, 1).Value = "wizard.new.import.database"
$ws.Cells.Item($insertRow, 2).Value = "Database"
$ws.Cells.Item($insertRow, 3).Value = "データベース"

# Grow the table (ListObject), and with it the AutoFilter range, so it
# spans through the last data row, which is now one row further down
# than before (C274 -> C275).
$lastRow = $lo.Range.Row + $lo.Range.Rows.Count
$firstCell = $ws.Cells.Item($lo.Range.Row, $lo.Range.Column)
$lastCell = $ws.Cells.Item($lastRow, $lo.Range.Column + $lo.Range.Columns.Count - 1)
$lo.Resize($ws.Range($firstCell, $lastCell))
